$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "56.876.53"
Set-TextValue "E2" "  +4.30%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.455.25"
Set-TextValue "E3" "  +1.86%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.27%  "

# Row 5 - BNB
Set-TextValue "D5" "490.68"
Set-TextValue "E5" "  +3.19%  "

# Row 6 - Solana
Set-TextValue "D6" "152.57"
Set-TextValue "E6" "  +10.61%  "

# Row 7 - USDC
Set-TextValue "E7" "  +0.11%  "

# Row 8 - XRP
Set-TextValue "E8" "  +2.72%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.456.41"
Set-TextValue "E9" "  +0.94%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.100"
Set-TextValue "E10" "  +5.43%  "

# Row 11 - Toncoin
Set-TextValue "D11" "5.70"
Set-TextValue "E11" "  +4.29%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.335"
Set-TextValue "E12" "  +3.92%  "

# Row 13 - TRON
Set-TextValue "D13" "0.125"
Set-TextValue "E13" "  +1.63%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "2.880.79"
Set-TextValue "E14" "  +1.69%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "57.106.73"
Set-TextValue "E15" "  +4.36%  "

# Row 16 - Avalanche
Set-TextValue "D16" "20.97"
Set-TextValue "E16" "  +3.30%  "

# Row 17 - ShibaInu
Set-TextValue "E17" "  +3.21%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.475.39"
Set-TextValue "E18" "  +1.63%  "

# Row 19 - Polkadot
Set-TextValue "D19" "4.59"
Set-TextValue "E19" "  +6.55%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "324.23"
Set-TextValue "E20" "  +4.40%  "

# Row 21 - Chainlink
Set-TextValue "E21" "  +1.73%  "

# Row 22 - Dai
Set-TextValue "E22" "  +0.33%  "

# Row 23 - Uniswap
Set-TextValue "E23" "  +4.02%  "

# Row 24 - Litecoin
Set-TextValue "D24" "58.10"
Set-TextValue "E24" "  +2.02%  "

# Row 25 - Polygon
Set-TextValue "E25" "  +2.03%  "

# Row 26 - Binance-PegBSC-USD
Set-TextValue "E26" "  -0.20%  "

# Row 27 - Kaspa
Set-TextValue "D27" "0.162"
Set-TextValue "E27" "  +1.32%  "

# Row 28 - WrappedeETH
Set-TextValue "D28" "2.567.70"
Set-TextValue "E28" "  +0.86%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue "E29" "  +4.08%  "

# Row 30 - PEPE
Set-TextValue "E30" "  +5.79%  "

# Row 31 - USDe
Set-TextValue "D31" "1.00"
Set-TextValue "E31" "  +0.21%  "

# Row 32 - Monero
Set-TextValue "D32" "151.05"
Set-TextValue "E32" "  +1.97%  "

# Row 33 - now PancakeSwap (was EthereumClassic)
Set-TextValue "B33" "PancakeSwap"
Set-TextValue "C33" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D33" "1.52"
Set-TextValue "E33" "  +4.40%  "

# Row 34 - now EthereumClassic (was PancakeSwap)
Set-TextValue "B34" "EthereumClassic"
Set-TextValue "C34" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D34" "18.21"
Set-TextValue "E34" "  +1.63%  "

# Row 35 - Aptos
Set-TextValue "D35" "5.24"
Set-TextValue "E35" "  +2.54%  "

# Row 36 - now Fetch.AI (was ImmutableX)
Set-TextValue "B36" "Fetch.AI"
Set-TextValue "C36" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D36" "0.890"
Set-TextValue "E36" "  +7.89%  "

# Row 37 - now ImmutableX (was NEARProtocol)
Set-TextValue "B37" "ImmutableX"
Set-TextValue "C37" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D37" "1.15"
Set-TextValue "E37" "  +3.51%  "

# Row 38 - now NEARProtocol (was Fetch.AI)
Set-TextValue "B38" "NEARProtocol"
Set-TextValue "C38" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D38" "3.75"
Set-TextValue "E38" "  +5.71%  "

# Row 39 - Stacks
Set-TextValue "D39" "1.40"
Set-TextValue "E39" "  +9.98%  "

# Row 40 - OKB
Set-TextValue "D40" "34.21"
Set-TextValue "E40" "  +2.27%  "

# Row 41 - Filecoin
Set-TextValue "D41" "3.51"
Set-TextValue "E41" "  +3.28%  "

# Row 42 - now Hedera (was FirstDigitalUSD)
Set-TextValue "B42" "Hedera"
Set-TextValue "C42" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D42" "0.0558"
Set-TextValue "E42" "  +3.13%  "

# Row 43 - now FirstDigitalUSD (was Hedera)
Set-TextValue "B43" "FirstDigitalUSD"
Set-TextValue "C43" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "0.995"
Set-TextValue "E43" "  +0.33%  "

# Row 44 - Mantle
Set-TextValue "D44" "0.606"
Set-TextValue "E44" "  +1.90%  "

# Row 45 - Stellar
Set-TextValue "D45" "0.0960"
Set-TextValue "E45" "  +7.75%  "

# Row 46 - now Bittensor (was RenderToken)
Set-TextValue "B46" "Bittensor"
Set-TextValue "C46" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D46" "266.15"
Set-TextValue "E46" "  +4.70%  "

# Row 47 - now RenderToken (was Bittensor)
Set-TextValue "B47" "RenderToken"
Set-TextValue "C47" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D47" "4.80"
Set-TextValue "E47" "  +2.09%  "

# Row 48 - WhiteBITCoin
Set-TextValue "E48" "  +0.85%  "

# Row 49 - VeChain
Set-TextValue "D49" "0.0229"
Set-TextValue "E49" "  +3.95%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "17.77"
Set-TextValue "E50" "  +4.61%  "

# Row 51 - dogwifhat
Set-TextValue "E51" "  +26.45%  "
